$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First names for the three registrants (amir / iravati / anu)
$ws.Range("A2").Value = "amir"
$ws.Range("A3").Value = "iravati"
$ws.Range("A4").Value = "anu"

# Last name column now reads "tester" for every registrant
$ws.Range("B2").Value = "tester"
$ws.Range("B3").Value = "tester"
$ws.Range("B4").Value = "tester"

# Telephone numbers stay the same digits, kept as text
$ws.Range("C2").Value = "'9267899098"
$ws.Range("C3").Value = "'9099909876"
$ws.Range("C4").Value = "'9898765432"

# Row 4's password changes to anu34
$ws.Range("D4").Value = "anu34"

# Move the active selection from C5 to A2
$ws.Range("A2").Select()
